$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 27: continue the work log with a new task entry.
# B27 = start time (19:40), C27 = end time (20:10), D27 = task description.
$ws.Range("B27").Value = 0.81944444444444453
$ws.Range("C27").Value = 0.84027777777777779
$ws.Range("B27:C27").NumberFormat = "h:mm"

$ws.Range("D27").Value = "Coding AvatarAnimationController (IK)"
$ws.Range("D27").WrapText = $true
$ws.Range("D27").HorizontalAlignment = -4131

# Update the selection to mirror where Excel would leave the cursor
# after entering this new row (next empty row in column D).
$ws.Range("D29").Select() | Out-Null
